$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(251442, 275600, 272396),
    @(236768, 275600, 272396),
    @(224416, 263344, 272396),
    @(224416, 257754, 272396),
    @(224416, 255380, 272396),
    @(224416, 256860, 272396),
    @(224416, 264924, 272396),
    @(224416, 263918, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396),
    @(224416, 275600, 272396)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}
